$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - first worksheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6460
$ws1.Range("F3").Value = 116
$ws1.Range("F5").Value = 389
$ws1.Range("F10").Value = 78
$ws1.Range("F12").Value = 157
$ws1.Range("F13").Value = 377
$ws1.Range("F14").Value = 947
$ws1.Range("F15").Value = 3167
$ws1.Range("F16").Value = 13
$ws1.Range("F17").Value = 193
$ws1.Range("F18").Value = 1840

# Sheet "全部类型" (All types) - fourth worksheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6460
$ws4.Range("F3").Value = 116
$ws4.Range("F5").Value = 389
$ws4.Range("F11").Value = 78
$ws4.Range("F13").Value = 157
$ws4.Range("F14").Value = 377
$ws4.Range("F15").Value = 947
$ws4.Range("F16").Value = 3167
$ws4.Range("F17").Value = 13
$ws4.Range("F18").Value = 193
$ws4.Range("F19").Value = 1840
